# Scheduled runner update: refresh market-price-derived columns (H-N)
# across several Leve profit tables, per the latest price snapshot.
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 3928.7334
$ws.Range("I70").Value = 1913.6666
$ws.Range("K70").Value = 5740.9998
$ws.Range("M70").Value = -5470.9998
$ws.Range("H73").Value = 3928.7334
$ws.Range("I73").Value = 1913.6666
$ws.Range("K73").Value = 5740.9998
$ws.Range("M73").Value = -4804.9998
$ws.Range("H107").Value = 2334
$ws.Range("I107").Value = 2334
$ws.Range("K107").Value = 2334
$ws.Range("M107").Value = -414
$ws.Range("H132").Value = 13748.55
$ws.Range("I132").Value = 12621.823
$ws.Range("J132").Value = 20133.334
$ws.Range("K132").Value = 37865.469
$ws.Range("L132").Value = 60400.00199999999
$ws.Range("M132").Value = -35335.469
$ws.Range("N132").Value = -65460.00199999999
$ws.Range("H135").Value = 893.5
$ws.Range("I135").Value = 894.6667
$ws.Range("K135").Value = 8052.0003
$ws.Range("M135").Value = -5517.0003
$ws.Range("H141").Value = 3266.6667
$ws.Range("I141").Value = 2125
$ws.Range("K141").Value = 6375
$ws.Range("M141").Value = -1195

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 4449.75
$ws.Range("I8").Value = 3725
$ws.Range("J8").Value = 5174.5
$ws.Range("K8").Value = 3725
$ws.Range("L8").Value = 5174.5
$ws.Range("M8").Value = -3581
$ws.Range("N8").Value = -5462.5
$ws.Range("H11").Value = 2616.75
$ws.Range("I11").Value = 1000
$ws.Range("J11").Value = 3155.6667
$ws.Range("K11").Value = 1000
$ws.Range("L11").Value = 3155.6667
$ws.Range("M11").Value = -856
$ws.Range("N11").Value = -3443.6667
$ws.Range("H13").Value = 4250
$ws.Range("J13").Value = 4250
$ws.Range("L13").Value = 4250
$ws.Range("N13").Value = -4538
$ws.Range("H74").Value = 5595.143
$ws.Range("I74").Value = 5630.1665
$ws.Range("K74").Value = 5630.1665
$ws.Range("M74").Value = -4756.1665
$ws.Range("H77").Value = 5595.143
$ws.Range("I77").Value = 5630.1665
$ws.Range("K77").Value = 28150.8325
$ws.Range("M77").Value = -23782.8325
$ws.Range("H88").Value = 1070
$ws.Range("I88").Value = 1506
$ws.Range("J88").Value = 721.2
$ws.Range("K88").Value = 1506
$ws.Range("L88").Value = 721.2
$ws.Range("M88").Value = -1100
$ws.Range("N88").Value = -1533.2
$ws.Range("H91").Value = 1070
$ws.Range("I91").Value = 1506
$ws.Range("J91").Value = 721.2
$ws.Range("K91").Value = 1506
$ws.Range("L91").Value = 721.2
$ws.Range("M91").Value = -102
$ws.Range("N91").Value = -3529.2
$ws.Range("H110").Value = 1042.7142
$ws.Range("I110").Value = 840.4167
$ws.Range("K110").Value = 840.4167
$ws.Range("M110").Value = 1204.5833

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 504.95834
$ws.Range("I6").Value = 504.95834
$ws.Range("K6").Value = 504.95834
$ws.Range("M6").Value = -391.95834
$ws.Range("H7").Value = 2743.077
$ws.Range("I7").Value = 4250.2915
$ws.Range("J7").Value = 331.53333
$ws.Range("K7").Value = 4250.2915
$ws.Range("L7").Value = 331.53333
$ws.Range("M7").Value = -4137.2915
$ws.Range("N7").Value = -557.53333
$ws.Range("H22").Value = 4398.8
$ws.Range("I22").Value = 4499.5
$ws.Range("K22").Value = 4499.5
$ws.Range("M22").Value = -4149.5
$ws.Range("H58").Value = 4622.727
$ws.Range("I58").Value = 3986.4285
$ws.Range("K58").Value = 3986.4285
$ws.Range("M58").Value = -3783.4285
$ws.Range("H132").Value = 1663.6666
$ws.Range("I132").Value = 1663.6666
$ws.Range("K132").Value = 4990.9998
$ws.Range("M132").Value = -2460.9998
$ws.Range("H134").Value = 2430.3076
$ws.Range("I134").Value = 2507.9167
$ws.Range("K134").Value = 7523.750100000001
$ws.Range("M134").Value = -4988.750100000001
$ws.Range("H136").Value = 4622.727
$ws.Range("I136").Value = 3986.4285
$ws.Range("K136").Value = 11959.2855
$ws.Range("M136").Value = -9409.2855

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 51515.82
$ws.Range("I4").Value = 64731.773
$ws.Range("J4").Value = 304
$ws.Range("K4").Value = 194195.319
$ws.Range("L4").Value = 912
$ws.Range("M4").Value = -194083.319
$ws.Range("N4").Value = -1136
$ws.Range("H23").Value = 275.83334
$ws.Range("I23").Value = 250
$ws.Range("J23").Value = 281
$ws.Range("K23").Value = 750
$ws.Range("L23").Value = 843
$ws.Range("M23").Value = -515
$ws.Range("N23").Value = -1313
$ws.Range("H37").Value = 188571.42
$ws.Range("J37").Value = 188571.42
$ws.Range("L37").Value = 565714.26
$ws.Range("N37").Value = -565938.26
$ws.Range("H55").Value = 5912.778
$ws.Range("J55").Value = 5912.778
$ws.Range("L55").Value = 17738.334
$ws.Range("N55").Value = -18092.334
$ws.Range("H61").Value = 198.57143
$ws.Range("I61").Value = 198.57143
$ws.Range("K61").Value = 595.71429
$ws.Range("M61").Value = -380.71429
$ws.Range("H112").Value = 913
$ws.Range("I112").Value = 913
$ws.Range("K112").Value = 2739
$ws.Range("M112").Value = -1631
$ws.Range("H115").Value = 864.3333
$ws.Range("I115").Value = 302.5
$ws.Range("J115").Value = 1988
$ws.Range("K115").Value = 907.5
$ws.Range("L115").Value = 5964
$ws.Range("M115").Value = 267.5
$ws.Range("N115").Value = -8314

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 466.66666
$ws.Range("I31").Value = 466.66666
$ws.Range("K31").Value = 466.66666
$ws.Range("M31").Value = -174.66666
$ws.Range("H37").Value = 466.66666
$ws.Range("I37").Value = 466.66666
$ws.Range("K37").Value = 466.66666
$ws.Range("M37").Value = -189.66666
$ws.Range("H132").Value = 3298.8
$ws.Range("I132").Value = 3298.8
$ws.Range("K132").Value = 9896.400000000001
$ws.Range("M132").Value = -7366.400000000001

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H132").Value = 4996
$ws.Range("I132").Value = 4996
$ws.Range("K132").Value = 14988
$ws.Range("M132").Value = -12458
$ws.Range("H136").Value = 2957.6667
$ws.Range("I136").Value = 2957.6667
$ws.Range("K136").Value = 8873.000100000001
$ws.Range("M136").Value = -6323.000100000001

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4750
$ws.Range("I81").Value = 2000
$ws.Range("K81").Value = 4000
$ws.Range("M81").Value = -2939
$ws.Range("H84").Value = 4750
$ws.Range("I84").Value = 2000
$ws.Range("K84").Value = 20000
$ws.Range("M84").Value = -14696
$ws.Range("H86").Value = 70324.336
$ws.Range("J86").Value = 70324.336
$ws.Range("L86").Value = 70324.336
$ws.Range("N86").Value = -72570.336
$ws.Range("H89").Value = 70324.336
$ws.Range("J89").Value = 70324.336
$ws.Range("L89").Value = 351621.68
$ws.Range("N89").Value = -362853.68
$ws.Range("H99").Value = 40000
$ws.Range("I99").Value = 40000
$ws.Range("K99").Value = 40000
$ws.Range("M99").Value = -37005
$ws.Range("H132").Value = 1939.1177
$ws.Range("I132").Value = 1712
$ws.Range("K132").Value = 5136
$ws.Range("M132").Value = -2606

